$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Single-cell corrections in existing rows (detect_structure / isPivot / backup fixes) ---
$ws.Cells.Item(65, 17).Value = 0   # Q65 (detect_structure): 1 -> 0
$ws.Cells.Item(72, 17).Value = 0   # Q72 (detect_structure): 1 -> 0
$ws.Cells.Item(75, 17).Value = 0   # Q75 (detect_structure): 2 -> 0
$ws.Cells.Item(735, 15).Value = 2  # O735 (isPivot): 0 -> 2
$ws.Cells.Item(737, 18).Value = 0  # R737 (backup): "" -> 0
$ws.Cells.Item(738, 18).Value = 0  # R738 (backup): "" -> 0

# --- Append 26 new weekly rows 739:764 (2024-07-01 .. 2024-12-23) ---
$dateFormat = $ws.Cells.Item(738, 1).NumberFormat

# row 739
$ws.Cells.Item(739, 1).Value = 45474
$ws.Cells.Item(739, 2).Value = 131.796933984319
$ws.Cells.Item(739, 3).Value = 143.4889989899596
$ws.Cells.Item(739, 4).Value = 129.866496992317
$ws.Cells.Item(739, 5).Value = 141.4690093994141
$ws.Cells.Item(739, 7).Value = 134803726
$ws.Cells.Item(739, 8).Value = 2024
$ws.Cells.Item(739, 9).Value = 7
$ws.Cells.Item(739, 10).Value = 1
$ws.Cells.Item(739, 11).Value = 0
$ws.Cells.Item(739, 12).Value = 0
$ws.Cells.Item(739, 13).Value = 0
$ws.Cells.Item(739, 14).Value = 27
$ws.Cells.Item(739, 15).Value = 0
$ws.Cells.Item(739, 16).Value = 0
$ws.Cells.Item(739, 17).Value = 2
$ws.Cells.Item(739, 1).NumberFormat = $dateFormat

# row 740
$ws.Cells.Item(740, 1).Value = 45481
$ws.Cells.Item(740, 2).Value = 142.5934441674954
$ws.Cells.Item(740, 3).Value = 156.2258915389616
$ws.Cells.Item(740, 4).Value = 135.3393799365587
$ws.Cells.Item(740, 5).Value = 149.9470062255859
$ws.Cells.Item(740, 7).Value = 250791130
$ws.Cells.Item(740, 8).Value = 2024
$ws.Cells.Item(740, 9).Value = 7
$ws.Cells.Item(740, 10).Value = 8
$ws.Cells.Item(740, 11).Value = 0
$ws.Cells.Item(740, 12).Value = 0
$ws.Cells.Item(740, 13).Value = 0
$ws.Cells.Item(740, 14).Value = 28
$ws.Cells.Item(740, 15).Value = 0
$ws.Cells.Item(740, 16).Value = 0
$ws.Cells.Item(740, 17).Value = 0
$ws.Cells.Item(740, 1).NumberFormat = $dateFormat

# row 741
$ws.Cells.Item(741, 1).Value = 45488
$ws.Cells.Item(741, 2).Value = 150.8027750194937
$ws.Cells.Item(741, 3).Value = 157.2209698005918
$ws.Cells.Item(741, 4).Value = 139.4888252856087
$ws.Cells.Item(741, 5).Value = 140.2052764892578
$ws.Cells.Item(741, 7).Value = 132487957
$ws.Cells.Item(741, 8).Value = 2024
$ws.Cells.Item(741, 9).Value = 7
$ws.Cells.Item(741, 10).Value = 15
$ws.Cells.Item(741, 11).Value = 0
$ws.Cells.Item(741, 12).Value = 0
$ws.Cells.Item(741, 13).Value = 0
$ws.Cells.Item(741, 14).Value = 29
$ws.Cells.Item(741, 15).Value = 0
$ws.Cells.Item(741, 16).Value = 0
$ws.Cells.Item(741, 17).Value = 1
$ws.Cells.Item(741, 1).NumberFormat = $dateFormat

# row 742
$ws.Cells.Item(742, 1).Value = 45495
$ws.Cells.Item(742, 2).Value = 139.2898111008772
$ws.Cells.Item(742, 3).Value = 158.8628244808739
$ws.Cells.Item(742, 4).Value = 129.4585357014445
$ws.Cells.Item(742, 5).Value = 147.3797302246094
$ws.Cells.Item(742, 7).Value = 176650499
$ws.Cells.Item(742, 8).Value = 2024
$ws.Cells.Item(742, 9).Value = 7
$ws.Cells.Item(742, 10).Value = 22
$ws.Cells.Item(742, 11).Value = 0
$ws.Cells.Item(742, 12).Value = 0
$ws.Cells.Item(742, 13).Value = 0
$ws.Cells.Item(742, 14).Value = 30
$ws.Cells.Item(742, 15).Value = 1
$ws.Cells.Item(742, 16).Value = 0
$ws.Cells.Item(742, 17).Value = 0
$ws.Cells.Item(742, 1).NumberFormat = $dateFormat

# row 743
$ws.Cells.Item(743, 1).Value = 45502
$ws.Cells.Item(743, 2).Value = 148.6633556858308
$ws.Cells.Item(743, 3).Value = 152.9421569525601
$ws.Cells.Item(743, 4).Value = 141.5585610926304
$ws.Cells.Item(743, 5).Value = 143.1009216308594
$ws.Cells.Item(743, 7).Value = 84379653
$ws.Cells.Item(743, 8).Value = 2024
$ws.Cells.Item(743, 9).Value = 7
$ws.Cells.Item(743, 10).Value = 29
$ws.Cells.Item(743, 11).Value = 0
$ws.Cells.Item(743, 12).Value = 0
$ws.Cells.Item(743, 13).Value = 0
$ws.Cells.Item(743, 14).Value = 31
$ws.Cells.Item(743, 15).Value = 0
$ws.Cells.Item(743, 16).Value = 0
$ws.Cells.Item(743, 17).Value = 0
$ws.Cells.Item(743, 1).NumberFormat = $dateFormat

# row 744
$ws.Cells.Item(744, 1).Value = 45509
$ws.Cells.Item(744, 2).Value = 139.3097117828249
$ws.Cells.Item(744, 3).Value = 143.0611275586574
$ws.Cells.Item(744, 4).Value = 134.3343649334383
$ws.Cells.Item(744, 5).Value = 138.6330718994141
$ws.Cells.Item(744, 7).Value = 65710387
$ws.Cells.Item(744, 8).Value = 2024
$ws.Cells.Item(744, 9).Value = 8
$ws.Cells.Item(744, 10).Value = 5
$ws.Cells.Item(744, 11).Value = 0
$ws.Cells.Item(744, 12).Value = 0
$ws.Cells.Item(744, 13).Value = 0
$ws.Cells.Item(744, 14).Value = 32
$ws.Cells.Item(744, 15).Value = 0
$ws.Cells.Item(744, 16).Value = 0
$ws.Cells.Item(744, 17).Value = 0
$ws.Cells.Item(744, 1).NumberFormat = $dateFormat

# row 745
$ws.Cells.Item(745, 1).Value = 45516
$ws.Cells.Item(745, 2).Value = 137.3195760883174
$ws.Cells.Item(745, 3).Value = 148.106124652941
$ws.Cells.Item(745, 4).Value = 134.9911173552206
$ws.Cells.Item(745, 5).Value = 142.255126953125
$ws.Cells.Item(745, 7).Value = 110047775
$ws.Cells.Item(745, 8).Value = 2024
$ws.Cells.Item(745, 9).Value = 8
$ws.Cells.Item(745, 10).Value = 12
$ws.Cells.Item(745, 11).Value = 0
$ws.Cells.Item(745, 12).Value = 0
$ws.Cells.Item(745, 13).Value = 0
$ws.Cells.Item(745, 14).Value = 33
$ws.Cells.Item(745, 15).Value = 0
$ws.Cells.Item(745, 16).Value = 0
$ws.Cells.Item(745, 17).Value = 0
$ws.Cells.Item(745, 1).NumberFormat = $dateFormat

# row 746
$ws.Cells.Item(746, 1).Value = 45523
$ws.Cells.Item(746, 2).Value = 143.0412255177244
$ws.Cells.Item(746, 3).Value = 143.0412255177244
$ws.Cells.Item(746, 4).Value = 131.2496471146678
$ws.Cells.Item(746, 5).Value = 133.0507354736328
$ws.Cells.Item(746, 7).Value = 66378363
$ws.Cells.Item(746, 8).Value = 2024
$ws.Cells.Item(746, 9).Value = 8
$ws.Cells.Item(746, 10).Value = 19
$ws.Cells.Item(746, 11).Value = 0
$ws.Cells.Item(746, 12).Value = 0
$ws.Cells.Item(746, 13).Value = 0
$ws.Cells.Item(746, 14).Value = 34
$ws.Cells.Item(746, 15).Value = 0
$ws.Cells.Item(746, 16).Value = 0
$ws.Cells.Item(746, 17).Value = 0
$ws.Cells.Item(746, 1).NumberFormat = $dateFormat

# row 747
$ws.Cells.Item(747, 1).Value = 45530
$ws.Cells.Item(747, 2).Value = 133.7373155612406
$ws.Cells.Item(747, 3).Value = 135.8269672852035
$ws.Cells.Item(747, 4).Value = 131.3491551769001
$ws.Cells.Item(747, 5).Value = 132.6925048828125
$ws.Cells.Item(747, 7).Value = 45788412
$ws.Cells.Item(747, 8).Value = 2024
$ws.Cells.Item(747, 9).Value = 8
$ws.Cells.Item(747, 10).Value = 26
$ws.Cells.Item(747, 11).Value = 0
$ws.Cells.Item(747, 12).Value = 0
$ws.Cells.Item(747, 13).Value = 0
$ws.Cells.Item(747, 14).Value = 35
$ws.Cells.Item(747, 15).Value = 0
$ws.Cells.Item(747, 16).Value = 0
$ws.Cells.Item(747, 17).Value = 0
$ws.Cells.Item(747, 1).NumberFormat = $dateFormat

# row 748
$ws.Cells.Item(748, 1).Value = 45537
$ws.Cells.Item(748, 2).Value = 139.3097148485176
$ws.Cells.Item(748, 3).Value = 139.3097148485176
$ws.Cells.Item(748, 4).Value = 129.1600100891221
$ws.Cells.Item(748, 5).Value = 129.8665008544922
$ws.Cells.Item(748, 7).Value = 50243329
$ws.Cells.Item(748, 8).Value = 2024
$ws.Cells.Item(748, 9).Value = 9
$ws.Cells.Item(748, 10).Value = 2
$ws.Cells.Item(748, 11).Value = 0
$ws.Cells.Item(748, 12).Value = 0
$ws.Cells.Item(748, 13).Value = 0
$ws.Cells.Item(748, 14).Value = 36
$ws.Cells.Item(748, 15).Value = 0
$ws.Cells.Item(748, 16).Value = 0
$ws.Cells.Item(748, 17).Value = 2
$ws.Cells.Item(748, 1).NumberFormat = $dateFormat

# row 749
$ws.Cells.Item(749, 1).Value = 45544
$ws.Cells.Item(749, 2).Value = 129.4585399045288
$ws.Cells.Item(749, 3).Value = 133.7771376196379
$ws.Cells.Item(749, 4).Value = 124.1448651600474
$ws.Cells.Item(749, 5).Value = 129.9262237548828
$ws.Cells.Item(749, 7).Value = 57120422
$ws.Cells.Item(749, 8).Value = 2024
$ws.Cells.Item(749, 9).Value = 9
$ws.Cells.Item(749, 10).Value = 9
$ws.Cells.Item(749, 11).Value = 0
$ws.Cells.Item(749, 12).Value = 0
$ws.Cells.Item(749, 13).Value = 0
$ws.Cells.Item(749, 14).Value = 37
$ws.Cells.Item(749, 15).Value = 0
$ws.Cells.Item(749, 16).Value = 0
$ws.Cells.Item(749, 17).Value = 0
$ws.Cells.Item(749, 1).NumberFormat = $dateFormat

# row 750
$ws.Cells.Item(750, 1).Value = 45551
$ws.Cells.Item(750, 2).Value = 130.8000030517578
$ws.Cells.Item(750, 3).Value = 132.3000030517578
$ws.Cells.Item(750, 4).Value = 125.5100021362305
$ws.Cells.Item(750, 5).Value = 130.3000030517578
$ws.Cells.Item(750, 7).Value = 32891735
$ws.Cells.Item(750, 8).Value = 2024
$ws.Cells.Item(750, 9).Value = 9
$ws.Cells.Item(750, 10).Value = 16
$ws.Cells.Item(750, 11).Value = 0
$ws.Cells.Item(750, 12).Value = 0
$ws.Cells.Item(750, 13).Value = 0
$ws.Cells.Item(750, 14).Value = 38
$ws.Cells.Item(750, 15).Value = 0
$ws.Cells.Item(750, 16).Value = 0
$ws.Cells.Item(750, 17).Value = 0
$ws.Cells.Item(750, 1).NumberFormat = $dateFormat

# row 751
$ws.Cells.Item(751, 1).Value = 45558
$ws.Cells.Item(751, 2).Value = 130.6499938964844
$ws.Cells.Item(751, 3).Value = 134.8899993896484
$ws.Cells.Item(751, 4).Value = 125.8000030517578
$ws.Cells.Item(751, 5).Value = 133.8399963378906
$ws.Cells.Item(751, 7).Value = 59604099
$ws.Cells.Item(751, 8).Value = 2024
$ws.Cells.Item(751, 9).Value = 9
$ws.Cells.Item(751, 10).Value = 23
$ws.Cells.Item(751, 11).Value = 0
$ws.Cells.Item(751, 12).Value = 0
$ws.Cells.Item(751, 13).Value = 0
$ws.Cells.Item(751, 14).Value = 39
$ws.Cells.Item(751, 15).Value = 0
$ws.Cells.Item(751, 16).Value = 0
$ws.Cells.Item(751, 17).Value = 0
$ws.Cells.Item(751, 1).NumberFormat = $dateFormat

# row 752
$ws.Cells.Item(752, 1).Value = 45565
$ws.Cells.Item(752, 2).Value = 132.9900054931641
$ws.Cells.Item(752, 3).Value = 135
$ws.Cells.Item(752, 4).Value = 124.0999984741211
$ws.Cells.Item(752, 5).Value = 126.6699981689453
$ws.Cells.Item(752, 7).Value = 36234299
$ws.Cells.Item(752, 8).Value = 2024
$ws.Cells.Item(752, 9).Value = 9
$ws.Cells.Item(752, 10).Value = 30
$ws.Cells.Item(752, 11).Value = 0
$ws.Cells.Item(752, 12).Value = 0
$ws.Cells.Item(752, 13).Value = 0
$ws.Cells.Item(752, 14).Value = 40
$ws.Cells.Item(752, 15).Value = 0
$ws.Cells.Item(752, 16).Value = 0
$ws.Cells.Item(752, 17).Value = 0
$ws.Cells.Item(752, 1).NumberFormat = $dateFormat

# row 753
$ws.Cells.Item(753, 1).Value = 45572
$ws.Cells.Item(753, 2).Value = 126.7099990844727
$ws.Cells.Item(753, 3).Value = 127.5
$ws.Cells.Item(753, 4).Value = 116.370002746582
$ws.Cells.Item(753, 5).Value = 120.1900024414062
$ws.Cells.Item(753, 7).Value = 40521070
$ws.Cells.Item(753, 8).Value = 2024
$ws.Cells.Item(753, 9).Value = 10
$ws.Cells.Item(753, 10).Value = 7
$ws.Cells.Item(753, 11).Value = 0
$ws.Cells.Item(753, 12).Value = 0
$ws.Cells.Item(753, 13).Value = 0
$ws.Cells.Item(753, 14).Value = 41
$ws.Cells.Item(753, 15).Value = 0
$ws.Cells.Item(753, 16).Value = 0
$ws.Cells.Item(753, 17).Value = 0
$ws.Cells.Item(753, 1).NumberFormat = $dateFormat

# row 754
$ws.Cells.Item(754, 1).Value = 45579
$ws.Cells.Item(754, 2).Value = 120.7099990844727
$ws.Cells.Item(754, 3).Value = 123.8600006103516
$ws.Cells.Item(754, 4).Value = 115.5199966430664
$ws.Cells.Item(754, 5).Value = 122.2600021362305
$ws.Cells.Item(754, 7).Value = 32408930
$ws.Cells.Item(754, 8).Value = 2024
$ws.Cells.Item(754, 9).Value = 10
$ws.Cells.Item(754, 10).Value = 14
$ws.Cells.Item(754, 11).Value = 0
$ws.Cells.Item(754, 12).Value = 0
$ws.Cells.Item(754, 13).Value = 0
$ws.Cells.Item(754, 14).Value = 42
$ws.Cells.Item(754, 15).Value = 0
$ws.Cells.Item(754, 16).Value = 0
$ws.Cells.Item(754, 17).Value = 0
$ws.Cells.Item(754, 1).NumberFormat = $dateFormat

# row 755
$ws.Cells.Item(755, 1).Value = 45586
$ws.Cells.Item(755, 2).Value = 122.9800033569336
$ws.Cells.Item(755, 3).Value = 122.9800033569336
$ws.Cells.Item(755, 4).Value = 106.1999969482422
$ws.Cells.Item(755, 5).Value = 107.2799987792969
$ws.Cells.Item(755, 7).Value = 36986332
$ws.Cells.Item(755, 8).Value = 2024
$ws.Cells.Item(755, 9).Value = 10
$ws.Cells.Item(755, 10).Value = 21
$ws.Cells.Item(755, 11).Value = 0
$ws.Cells.Item(755, 12).Value = 0
$ws.Cells.Item(755, 13).Value = 0
$ws.Cells.Item(755, 14).Value = 43
$ws.Cells.Item(755, 15).Value = 0
$ws.Cells.Item(755, 16).Value = 0
$ws.Cells.Item(755, 17).Value = 0
$ws.Cells.Item(755, 1).NumberFormat = $dateFormat

# row 756
$ws.Cells.Item(756, 1).Value = 45593
$ws.Cells.Item(756, 2).Value = 107.9800033569336
$ws.Cells.Item(756, 3).Value = 115.0999984741211
$ws.Cells.Item(756, 4).Value = 106.2200012207031
$ws.Cells.Item(756, 5).Value = 114.9000015258789
$ws.Cells.Item(756, 7).Value = 19118002
$ws.Cells.Item(756, 8).Value = 2024
$ws.Cells.Item(756, 9).Value = 10
$ws.Cells.Item(756, 10).Value = 28
$ws.Cells.Item(756, 11).Value = 0
$ws.Cells.Item(756, 12).Value = 0
$ws.Cells.Item(756, 13).Value = 0
$ws.Cells.Item(756, 14).Value = 44
$ws.Cells.Item(756, 15).Value = 0
$ws.Cells.Item(756, 16).Value = 0
$ws.Cells.Item(756, 17).Value = 0
$ws.Cells.Item(756, 1).NumberFormat = $dateFormat

# row 757
$ws.Cells.Item(757, 1).Value = 45600
$ws.Cells.Item(757, 2).Value = 114.9899978637695
$ws.Cells.Item(757, 3).Value = 117.4000015258789
$ws.Cells.Item(757, 4).Value = 109.1500015258789
$ws.Cells.Item(757, 5).Value = 111.4599990844727
$ws.Cells.Item(757, 7).Value = 31703150
$ws.Cells.Item(757, 8).Value = 2024
$ws.Cells.Item(757, 9).Value = 11
$ws.Cells.Item(757, 10).Value = 4
$ws.Cells.Item(757, 11).Value = 0
$ws.Cells.Item(757, 12).Value = 0
$ws.Cells.Item(757, 13).Value = 0
$ws.Cells.Item(757, 14).Value = 45
$ws.Cells.Item(757, 15).Value = 0
$ws.Cells.Item(757, 16).Value = 0
$ws.Cells.Item(757, 17).Value = 0
$ws.Cells.Item(757, 1).NumberFormat = $dateFormat

# row 758
$ws.Cells.Item(758, 1).Value = 45607
$ws.Cells.Item(758, 2).Value = 110.0199966430664
$ws.Cells.Item(758, 3).Value = 111.7699966430664
$ws.Cells.Item(758, 4).Value = 103.3600006103516
$ws.Cells.Item(758, 5).Value = 103.7699966430664
$ws.Cells.Item(758, 7).Value = 20448078
$ws.Cells.Item(758, 8).Value = 2024
$ws.Cells.Item(758, 9).Value = 11
$ws.Cells.Item(758, 10).Value = 11
$ws.Cells.Item(758, 11).Value = 0
$ws.Cells.Item(758, 12).Value = 0
$ws.Cells.Item(758, 13).Value = 0
$ws.Cells.Item(758, 14).Value = 46
$ws.Cells.Item(758, 15).Value = 0
$ws.Cells.Item(758, 16).Value = 0
$ws.Cells.Item(758, 17).Value = 0
$ws.Cells.Item(758, 1).NumberFormat = $dateFormat

# row 759
$ws.Cells.Item(759, 1).Value = 45614
$ws.Cells.Item(759, 2).Value = 104
$ws.Cells.Item(759, 3).Value = 110.3000030517578
$ws.Cells.Item(759, 4).Value = 101.6999969482422
$ws.Cells.Item(759, 5).Value = 107.5699996948242
$ws.Cells.Item(759, 7).Value = 43343378
$ws.Cells.Item(759, 8).Value = 2024
$ws.Cells.Item(759, 9).Value = 11
$ws.Cells.Item(759, 10).Value = 18
$ws.Cells.Item(759, 11).Value = 0
$ws.Cells.Item(759, 12).Value = 0
$ws.Cells.Item(759, 13).Value = 0
$ws.Cells.Item(759, 14).Value = 47
$ws.Cells.Item(759, 15).Value = 2
$ws.Cells.Item(759, 16).Value = 0
$ws.Cells.Item(759, 17).Value = 0
$ws.Cells.Item(759, 1).NumberFormat = $dateFormat

# row 760
$ws.Cells.Item(760, 1).Value = 45621
$ws.Cells.Item(760, 2).Value = 109.8499984741211
$ws.Cells.Item(760, 3).Value = 119.6999969482422
$ws.Cells.Item(760, 4).Value = 109.75
$ws.Cells.Item(760, 5).Value = 114.8899993896484
$ws.Cells.Item(760, 7).Value = 37249132
$ws.Cells.Item(760, 8).Value = 2024
$ws.Cells.Item(760, 9).Value = 11
$ws.Cells.Item(760, 10).Value = 25
$ws.Cells.Item(760, 11).Value = 0
$ws.Cells.Item(760, 12).Value = 0
$ws.Cells.Item(760, 13).Value = 0
$ws.Cells.Item(760, 14).Value = 48
$ws.Cells.Item(760, 15).Value = 0
$ws.Cells.Item(760, 16).Value = 0
$ws.Cells.Item(760, 17).Value = 0
$ws.Cells.Item(760, 1).NumberFormat = $dateFormat

# row 761
$ws.Cells.Item(761, 1).Value = 45628
$ws.Cells.Item(761, 2).Value = 114.5
$ws.Cells.Item(761, 3).Value = 122.8000030517578
$ws.Cells.Item(761, 4).Value = 113.6100006103516
$ws.Cells.Item(761, 5).Value = 120.6100006103516
$ws.Cells.Item(761, 7).Value = 40223735
$ws.Cells.Item(761, 8).Value = 2024
$ws.Cells.Item(761, 9).Value = 12
$ws.Cells.Item(761, 10).Value = 2
$ws.Cells.Item(761, 11).Value = 0
$ws.Cells.Item(761, 12).Value = 0
$ws.Cells.Item(761, 13).Value = 0
$ws.Cells.Item(761, 14).Value = 49
$ws.Cells.Item(761, 15).Value = 0
$ws.Cells.Item(761, 16).Value = 0
$ws.Cells.Item(761, 17).Value = 0
$ws.Cells.Item(761, 1).NumberFormat = $dateFormat

# row 762
$ws.Cells.Item(762, 1).Value = 45635
$ws.Cells.Item(762, 2).Value = 121.1500015258789
$ws.Cells.Item(762, 3).Value = 124.5
$ws.Cells.Item(762, 4).Value = 115
$ws.Cells.Item(762, 5).Value = 118.0599975585938
$ws.Cells.Item(762, 7).Value = 33246083
$ws.Cells.Item(762, 8).Value = 2024
$ws.Cells.Item(762, 9).Value = 12
$ws.Cells.Item(762, 10).Value = 9
$ws.Cells.Item(762, 11).Value = 0
$ws.Cells.Item(762, 12).Value = 0
$ws.Cells.Item(762, 13).Value = 0
$ws.Cells.Item(762, 14).Value = 50
$ws.Cells.Item(762, 15).Value = 0
$ws.Cells.Item(762, 16).Value = 0
$ws.Cells.Item(762, 17).Value = 0
$ws.Cells.Item(762, 1).NumberFormat = $dateFormat

# row 763
$ws.Cells.Item(763, 1).Value = 45642
$ws.Cells.Item(763, 2).Value = 118.0599975585938
$ws.Cells.Item(763, 3).Value = 119.9499969482422
$ws.Cells.Item(763, 4).Value = 111
$ws.Cells.Item(763, 5).Value = 111.5400009155273
$ws.Cells.Item(763, 7).Value = 24828209
$ws.Cells.Item(763, 8).Value = 2024
$ws.Cells.Item(763, 9).Value = 12
$ws.Cells.Item(763, 10).Value = 16
$ws.Cells.Item(763, 11).Value = 0
$ws.Cells.Item(763, 12).Value = 0
$ws.Cells.Item(763, 13).Value = 0
$ws.Cells.Item(763, 14).Value = 51
$ws.Cells.Item(763, 15).Value = 0
$ws.Cells.Item(763, 16).Value = 0
$ws.Cells.Item(763, 17).Value = 0
$ws.Cells.Item(763, 1).NumberFormat = $dateFormat

# row 764
$ws.Cells.Item(764, 1).Value = 45649
$ws.Cells.Item(764, 2).Value = 112.2900009155273
$ws.Cells.Item(764, 3).Value = 112.7200012207031
$ws.Cells.Item(764, 4).Value = 107
$ws.Cells.Item(764, 5).Value = 107.4599990844727
$ws.Cells.Item(764, 7).Value = 15519435
$ws.Cells.Item(764, 8).Value = 2024
$ws.Cells.Item(764, 9).Value = 12
$ws.Cells.Item(764, 10).Value = 23
$ws.Cells.Item(764, 11).Value = 0
$ws.Cells.Item(764, 12).Value = 0
$ws.Cells.Item(764, 13).Value = 0
$ws.Cells.Item(764, 14).Value = 52
$ws.Cells.Item(764, 15).Value = 0
$ws.Cells.Item(764, 16).Value = 0
$ws.Cells.Item(764, 17).Value = 0
$ws.Cells.Item(764, 1).NumberFormat = $dateFormat
